$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 05:43"

# --- Countries re-sorted by total cases: Kazajistan overtakes Japon ---
# Row 54 now shows Kazajistan with freshly updated totals; row 55 keeps
# the country label "Japon" (which now holds what used to be Japon's data).
$ws.Range("A54").Value = "Kazajistan"
$ws.Range("B54").Value = 18231
$ws.Range("C54").Value = 499
$ws.Range("D54").Value = 11158
$ws.Range("E54").Value = 6946
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 127

$ws.Range("A55").Value = "Japon"
$ws.Range("B55").Value = 17916
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 16133
$ws.Range("E55").Value = 830
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 953

# --- Countries re-sorted by total cases: Honduras overtakes Azerbaiyan/Serbia ---
$ws.Range("A60").Value = "Honduras"
$ws.Range("B60").Value = 13356
$ws.Range("C60").Value = 584
$ws.Range("D60").Value = 1362
$ws.Range("E60").Value = 11599
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 32
$ws.Range("H60").Value = 395

$ws.Range("A61").Value = "Azerbaiyan"
$ws.Range("B61").Value = 13207
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 7168
$ws.Range("E61").Value = 5878
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 161

$ws.Range("A62").Value = "Serbia"
$ws.Range("B62").Value = 12990
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 11997
$ws.Range("E62").Value = 731
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 262

# --- Australia (row 74) data refresh ---
$ws.Range("B74").Value = 7492
$ws.Range("C74").Value = 18
$ws.Range("D74").Value = 6904
$ws.Range("E74").Value = 486

# --- Surinam (row 160) data refresh ---
$ws.Range("D160").Value = 132
$ws.Range("E160").Value = 179

# --- Birmania (row 161) data refresh ---
$ws.Range("B161").Value = 291
$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 85

# --- Mongolia (row 165) data refresh ---
$ws.Range("B165").Value = 215
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 158
$ws.Range("E165").Value = 57
